$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# C1 / D1 are brand new header cells: give them the exact same style as
# the existing header cells (bold font, thin border, centered/top aligned)
# by copying the formatting straight from A1, rather than rebuilding it
# property-by-property (which would mint a near-duplicate style).
$ws.Range("A1").Copy($ws.Range("C1"))
$ws.Range("A1").Copy($ws.Range("D1"))

$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "Value_1"
$ws.Range("C1").Value = "Value_2"
$ws.Range("D1").Value = "Value_3"

# --- Data rows: Label, Col B, Col C, Col D ---
$data = @(
    @("Revenues", "282,836", "307,394", "350,018"),
    @("Cost of revenues", "126,203", "133,332", "146,306"),
    @("Research and development", "39,500", "45,427", "49,326"),
    @("Sales and marketing", "26,567", "27,917", "27,808"),
    @("General and administrative", "15,724", "16,425", "14,188"),
    @("Total costs and expenses", "207,994", "223,101", "237,628"),
    @("Income from operations", "74,842", "84,293", "112,390"),
    @("Other income (expense), net (3,514)", "1,424", "7,425", $null),
    @("Income before income taxes", "71,328", "85,717", "119,815"),
    @("Provision for income taxes", "11,356", "11,922", "19,697"),
    @("Net income", "59,972", "73,795", "100,118")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    # Columns B/C/D hold numbers-with-commas that must stay *text*.
    # Force text entry via NumberFormat="@", then clear the format again
    # so the cell keeps its default (no explicit) style like the rest of
    # the sheet while remaining a text value.
    for ($col = 2; $col -le 4; $col++) {
        $text = $row[$col - 1]
        if ($text -ne $null) {
            $cell = $ws.Cells.Item($r, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $text
            $cell.ClearFormats()
        }
    }
    $r++
}

# Remove rows 13-15 (old extra rows beyond the new 12-row table)
$ws.Range("A13:D15").Delete() | Out-Null
